$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new rows first so all row numbers land in their final
# positions before any cell content is written.
#
# Row 18 (new): "igm_iga_result" - inserted just after "igm_igg_result"
# (old row 17), pushing "control" and everything below it down by one.
$ws.Rows("18:18").Insert()

# Row 31 (new, after the first insert has shifted things down): "igm_iga_agree"
# - inserted just before the final "antibody_agree" row, pushing it to row 32.
$ws.Rows("31:31").Insert()

# Fill in the new "igm_iga_agree" row content.
$ws.Range("B31").Value = "igm_iga_agree"
$ws.Range("D31").Value = "Agreement between igm_iga_result and antibody_truth."
$ws.Range("C31").Value = "string"
$ws.Rows("31:31").RowHeight = 17

# Fill in the new "igm_iga_result" row content.
$ws.Range("B18").Value = "igm_iga_result"
$ws.Range("D18").Value = "The test result for qualitative detection of (IgM / IgA) combined antibodies."
$ws.Range("C18").Value = "string"
$ws.Rows("18:18").RowHeight = 17

# Match the saved view state (scrolled down, D19 selected).
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()
